$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 45. This shifts the existing rows 45-58 down to 46-59,
# preserving all of their data untouched.
$ws.Rows("45").Insert()

# Populate the newly inserted row 45 with this week's new data record
# (same market/category/variety/quality/unit/origin as the prior top row,
# new date + new price figures).
$ws.Cells.Item(45, 1).Value = 11
$ws.Cells.Item(45, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(45, 3).Value = "Bíobío"
$ws.Cells.Item(45, 4).Value = 44825
$ws.Cells.Item(45, 5).Value = 8
$ws.Cells.Item(45, 6).Value = 100112031
$ws.Cells.Item(45, 7).Value = "Poroto verde"
$ws.Cells.Item(45, 8).Value = "Magnum"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 100
$ws.Cells.Item(45, 11).Value = 30000
$ws.Cells.Item(45, 12).Value = 32000
$ws.Cells.Item(45, 13).Value = 31000
$ws.Cells.Item(45, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(45, 15).Value = "Perú"
$ws.Cells.Item(45, 16).Value = 1240
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = "Hortaliza"
